$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "We can" + [char]8217 + "t wait to meet you! ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Kami tidak sabar untuk bertemu dengan Anda! ", 2)

$d.Content.Find.Execute(
    "In this email, we" + [char]8217 + "ve linked/attached the following documents:", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dalam email ini, kami telah melampirkan beberapa dokumen sebagai berikut:", 2)

$d.Content.Find.Execute(
    "Your return flight tickets", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Tiket penerbangan Anda", 2)

$d.Content.Find.Execute(
    "Your accommodation booking details", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Detail pemesanan akomodasi Anda", 2)
